$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a numeric-looking string as literal text,
    # matching the inlineStr/shared-string cells already in the sheet.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "27.823.46"
$ws.Range("E2").Value = "  +1.80%  "

# Row 3
$ws.Range("D3").Value = "1.887.57"
$ws.Range("E3").Value = "  +1.80%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.008"
$ws.Range("E4").Value = "  +0.51%  "

# Row 5
Set-TextValue $ws.Range("D5") "334.81"
$ws.Range("E5").Value = "  +1.77%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.008"
$ws.Range("E6").Value = "  +0.54%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4724"
$ws.Range("E7").Value = "  +2.15%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3930"
$ws.Range("E8").Value = "  -0.32%  "

# Row 9
Set-TextValue $ws.Range("D9") "47.55"
$ws.Range("E9").Value = "  +1.42%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.08073"
$ws.Range("E10").Value = "  +1.61%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.028"
$ws.Range("E11").Value = "  +1.52%  "

# Row 12
Set-TextValue $ws.Range("D12") "22.10"
$ws.Range("E12").Value = "  +2.90%  "

# Row 13
$ws.Range("D13").Value = "1.884.32"
$ws.Range("E13").Value = "  +3.51%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.993"
$ws.Range("E14").Value = "  +1.12%  "

# Row 15
Set-TextValue $ws.Range("D15") "7.140"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
Set-TextValue $ws.Range("D16") "1.011"
$ws.Range("E16").Value = "  +0.80%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.06749"
$ws.Range("E17").Value = "  +3.14%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.00001052"
$ws.Range("E18").Value = "  +2.08%  "

# Row 19
Set-TextValue $ws.Range("D19") "87.28"
$ws.Range("E19").Value = "  +1.30%  "

# Row 20
Set-TextValue $ws.Range("D20") "17.33"
$ws.Range("E20").Value = "  +0.95%  "

# Row 21
$ws.Range("E21").Value = "  +0.52%  "

# Row 22
$ws.Range("D22").Value = "27.854.42"
$ws.Range("E22").Value = "  +1.90%  "

# Row 23
Set-TextValue $ws.Range("D23") "5.526"
$ws.Range("E23").Value = "  +0.96%  "

# Row 24
Set-TextValue $ws.Range("D24") "11.00"
$ws.Range("E24").Value = "  +1.01%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.333"
$ws.Range("E25").Value = "  +1.31%  "

# Row 26
$ws.Range("D26").Value = "2.102.78"
$ws.Range("E26").Value = "  +2.56%  "

# Row 27
Set-TextValue $ws.Range("D27") "159.19"
$ws.Range("E27").Value = "  +3.71%  "

# Row 28
Set-TextValue $ws.Range("D28") "20.16"
$ws.Range("E28").Value = "  -1.40%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.104"
$ws.Range("E29").Value = "  +1.99%  "

# Row 30
Set-TextValue $ws.Range("D30") "5.577"
$ws.Range("E30").Value = "  +2.04%  "

# Row 31
Set-TextValue $ws.Range("D31") "122.05"
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.9794"

# Row 33
Set-TextValue $ws.Range("D33") "0.09506"
$ws.Range("E33").Value = "  +0.59%  "

# Row 34
$ws.Range("E34").Value = "  +1.20%  "

# Row 35
Set-TextValue $ws.Range("D35") "3.625"
$ws.Range("E35").Value = "  +1.09%  "

# Row 36
Set-TextValue $ws.Range("D36") "5.359"
$ws.Range("E36").Value = "  +1.93%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.06161"
$ws.Range("E37").Value = "  +1.97%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.02271"
$ws.Range("E38").Value = "  +1.98%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.221"
$ws.Range("E39").Value = "  +0.78%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D40") "0.6015"
$ws.Range("E40").Value = "  +1.46%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "8.057"
$ws.Range("E41").Value = "  +0.44%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.1898"
$ws.Range("E42").Value = "  +0.36%  "

# Row 43
Set-TextValue $ws.Range("D43") "10.31"
$ws.Range("E43").Value = "  +1.17%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.261"
$ws.Range("E44").Value = "  -1.60%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.5712"
$ws.Range("E45").Value = "  +1.64%  "

# Row 46
Set-TextValue $ws.Range("D46") "12.19"
$ws.Range("E46").Value = "  +0.49%  "

# Row 47
$ws.Range("E47").Value = "  -0.78%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.943"
$ws.Range("E48").Value = "  +1.15%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.06914"
$ws.Range("E49").Value = "  +2.31%  "

# Row 50
Set-TextValue $ws.Range("D50") "113.43"
$ws.Range("E50").Value = "  +3.70%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.00000000302"
$ws.Range("E51").Value = "  +8.13%  "
